# Apply the weekly cryptos-list refresh (GitHub Actions data pull).
# Cells hold plain text (prices/links/percent strings), not numeric values,
# so force the Text number format before writing to avoid Excel
# auto-coercing price strings like "398.10" into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.166.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.43%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.262.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.38%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "398.10"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.05%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.04"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.83%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.580"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.47%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.22%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.21%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0957"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.62%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.73%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.773.12"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.29%  "

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.46%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.00"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.53%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.264.36"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.28%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.03"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.79%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.05"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.58%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "56.992.37"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.33%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.20%  "

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.23%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.98"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.91%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "294.48"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.33%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.22"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.20%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.40%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.13"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.44%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.94"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.22%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.78%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.11%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.169"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.89%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.03%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.45%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.78%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.23"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +11.31%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0490"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.40%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.19%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.32"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.12%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.05%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.51%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.15%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "137.14"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.99%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.60%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.56%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.284"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.45%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.88"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.59%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.84"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.75%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.47"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.94%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.50%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.149.60"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.54%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -6.24%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.61%  "
